$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Default with Columns (BM25 0.7-0.3)", 0.4474358974358975, 0.6, 52.1882152557373),
    @("Default with Columns (Balanced 0.5-0.5)", 0.5487179487179488, 0.8153846153846154, 54.05382180213928),
    @("Default with Columns (Vector 0.3-0.7)", 0.6205128205128204, 0.8307692307692308, 55.41543316841125),
    @("Advanced with Columns (No Rerank)", 0.4671794871794872, 0.5692307692307692, 61.1391978263855),
    @("Advanced with Columns (With Rerank)", 0.6735897435897436, 0.8461538461538461, 54.67233204841614)
)

$startRow = 15
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
